# Commit: "updated version to 3.4"
#
# The only substantive content change in the source diff is the template's
# version-number textbox, which lives on the "Welkom" custom (slide) layout
# used by the deck's master -- not on an actual slide -- so it must be
# reached via the SlideMaster's CustomLayouts collection.
#
# All the other hunks in the diff are pure xmlns attribute re-orderings on
# unrelated, content-identical extension blocks (p14:creationId,
# a14:useLocalDpi, a14:hiddenFill/hiddenLine) introduced incidentally by the
# upstream tool that produced the diff; they carry no semantic change and
# are not reachable/expressible via the PowerPoint object model, so they are
# intentionally left untouched here.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

$found = $false
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shape = $layout.Shapes.Item($si)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -eq "3.3.4") {
                $shape.TextFrame.TextRange.Text = "3.4"
                $found = $true
            }
        }
    }
}

Write-Output "updated: $found"
